$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fndc5"
$ws.Range("C2").Value = "Itgav"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.116649
$ws.Range("H2").Value = 0.349947
$ws.Range("I2").Value = 0.186245347817965
$ws.Range("J2").Value = 0.186245347817965
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 2.461217067192
$ws.Range("R2").Value = 22.150953604728
$ws.Range("S2").Value = 0.05433216715089719
$ws.Range("T2").Value = 0.0543321671508972

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fndc5"
$ws.Range("C3").Value = "Itgav"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.116649
$ws.Range("H3").Value = 0.349947
$ws.Range("I3").Value = 0.186245347817965
$ws.Range("J3").Value = 0.186245347817965
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 4.178301584379
$ws.Range("R3").Value = 37.604714259411
$ws.Range("S3").Value = 0.09223736626706433
$ws.Range("T3").Value = 0.09223736626706434

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fndc5"
$ws.Range("C4").Value = "Itgav"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.116649
$ws.Range("H4").Value = 0.349947
$ws.Range("I4").Value = 0.186245347817965
$ws.Range("J4").Value = 0.186245347817965
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 1.79729240847
$ws.Range("R4").Value = 16.17563167623
$ws.Range("S4").Value = 0.03967581440000346
$ws.Range("T4").Value = 0.03967581440000346

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fndc5"
$ws.Range("C5").Value = "Itgav"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.50967
$ws.Range("H5").Value = 1.52901
$ws.Range("I5").Value = 0.813754652182035
$ws.Range("J5").Value = 0.813754652182035
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 10.75370129736
$ws.Range("R5").Value = 96.78331167624
$ws.Range("S5").Value = 0.2373914532640466
$ws.Range("T5").Value = 0.2373914532640466

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fndc5"
$ws.Range("C6").Value = "Itgav"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.50967
$ws.Range("H6").Value = 1.52901
$ws.Range("I6").Value = 0.813754652182035
$ws.Range("J6").Value = 0.813754652182035
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 18.25609279557
$ws.Range("R6").Value = 164.30483516013
$ws.Range("S6").Value = 0.4030091853795118
$ws.Range("T6").Value = 0.4030091853795119

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fndc5"
$ws.Range("C7").Value = "Itgav"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.50967
$ws.Range("H7").Value = 1.52901
$ws.Range("I7").Value = 0.813754652182035
$ws.Range("J7").Value = 0.813754652182035
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 7.852840760099999
$ws.Range("R7").Value = 70.6755668409
$ws.Range("S7").Value = 0.1733540135384766
$ws.Range("T7").Value = 0.1733540135384767

